$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.441.49"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.474.67"
$ws.Range("E3").Value = "  -3.00%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'587.25"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "'172.39"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("D9").Value = "2.472.21"
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D12").Value = "'5.09"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "'0.339"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").Value = "'26.14"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D16").Value = "'0.0000176"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "67.418.53"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "2.455.29"
$ws.Range("E18").Value = "  -4.00%  "
$ws.Range("D19").Value = "'11.68"
$ws.Range("E19").Value = "  +3.12%  "
$ws.Range("D20").Value = "'7.88"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "'366.04"
$ws.Range("E21").Value = "  +2.89%  "
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'71.11"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "2.621.21"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "'529.73"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "'1.30"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -3.57%  "
$ws.Range("D37").Value = "'159.19"
$ws.Range("D38").Value = "'1.41"
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("D39").Value = "'18.66"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").Value = "'18.61"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.09"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.76"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "'1.01"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "'2.47"
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0277"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'144.17"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("E51").Value = "  -2.13%  "
